$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 45033
$ws.Cells.Item(2, 15).Value = 16000
$ws.Cells.Item(2, 16).Value = 15500
$ws.Cells.Item(2, 19).Value = 861

# Row 3
$ws.Cells.Item(3, 4).Value = 45014
$ws.Cells.Item(3, 13).Value = 30
$ws.Cells.Item(3, 14).Value = 18000
$ws.Cells.Item(3, 15).Value = 18000
$ws.Cells.Item(3, 16).Value = 18000
$ws.Cells.Item(3, 19).Value = 1000

# Row 4
$ws.Cells.Item(4, 4).Value = 45091
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 22000
$ws.Cells.Item(4, 15).Value = 22000
$ws.Cells.Item(4, 16).Value = 22000
$ws.Cells.Item(4, 19).Value = 1222

# Row 5
$ws.Cells.Item(5, 4).Value = 45089
$ws.Cells.Item(5, 14).Value = 22000
$ws.Cells.Item(5, 15).Value = 23000
$ws.Cells.Item(5, 16).Value = 22500
$ws.Cells.Item(5, 19).Value = 1250

# Row 6
$ws.Cells.Item(6, 4).Value = 45037
$ws.Cells.Item(6, 14).Value = 16000
$ws.Cells.Item(6, 15).Value = 16000
$ws.Cells.Item(6, 16).Value = 16000
$ws.Cells.Item(6, 19).Value = 889

# Row 7
$ws.Cells.Item(7, 4).Value = 45062
$ws.Cells.Item(7, 13).Value = 90
$ws.Cells.Item(7, 14).Value = 13000
$ws.Cells.Item(7, 15).Value = 14000
$ws.Cells.Item(7, 16).Value = 13444
$ws.Cells.Item(7, 19).Value = 747

# Row 8
$ws.Cells.Item(8, 4).Value = 44999
$ws.Cells.Item(8, 13).Value = 60
$ws.Cells.Item(8, 14).Value = 17000
$ws.Cells.Item(8, 16).Value = 17500
$ws.Cells.Item(8, 19).Value = 972

# Row 9
$ws.Cells.Item(9, 4).Value = 45049
$ws.Cells.Item(9, 13).Value = 80
$ws.Cells.Item(9, 14).Value = 15000
$ws.Cells.Item(9, 15).Value = 15000
$ws.Cells.Item(9, 16).Value = 15000
$ws.Cells.Item(9, 19).Value = 833

# Row 10
$ws.Cells.Item(10, 4).Value = 45020
$ws.Cells.Item(10, 13).Value = 50
$ws.Cells.Item(10, 14).Value = 15000
$ws.Cells.Item(10, 15).Value = 15000
$ws.Cells.Item(10, 16).Value = 15000
$ws.Cells.Item(10, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(10, 19).Value = 938
$ws.Cells.Item(10, 20).Value = 16

# Row 11
$ws.Cells.Item(11, 4).Value = 45050
$ws.Cells.Item(11, 13).Value = 40
$ws.Cells.Item(11, 14).Value = 14000
$ws.Cells.Item(11, 15).Value = 14000
$ws.Cells.Item(11, 16).Value = 14000
$ws.Cells.Item(11, 19).Value = 778

# Row 12
$ws.Cells.Item(12, 4).Value = 45028
$ws.Cells.Item(12, 13).Value = 50
$ws.Cells.Item(12, 14).Value = 18000
$ws.Cells.Item(12, 16).Value = 18000
$ws.Cells.Item(12, 19).Value = 1000

# Row 13
$ws.Cells.Item(13, 4).Value = 45044
$ws.Cells.Item(13, 13).Value = 60
$ws.Cells.Item(13, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(13, 18).Value = "Región Metropolitana"
$ws.Cells.Item(13, 19).Value = 833
$ws.Cells.Item(13, 20).Value = 18

# Row 14
$ws.Cells.Item(14, 4).Value = 45002
$ws.Cells.Item(14, 13).Value = 30
$ws.Cells.Item(14, 14).Value = 18000
$ws.Cells.Item(14, 15).Value = 18000
$ws.Cells.Item(14, 16).Value = 18000
$ws.Cells.Item(14, 19).Value = 1000

# Row 15
$ws.Cells.Item(15, 4).Value = 45043
$ws.Cells.Item(15, 13).Value = 60
$ws.Cells.Item(15, 14).Value = 15000
$ws.Cells.Item(15, 15).Value = 15000
$ws.Cells.Item(15, 16).Value = 15000
$ws.Cells.Item(15, 19).Value = 833

# Row 16
$ws.Cells.Item(16, 4).Value = 45041
$ws.Cells.Item(16, 15).Value = 15000
$ws.Cells.Item(16, 16).Value = 15000
$ws.Cells.Item(16, 19).Value = 833

# Row 19
$ws.Cells.Item(19, 4).Value = 45030
$ws.Cells.Item(19, 13).Value = 40
$ws.Cells.Item(19, 14).Value = 18000
$ws.Cells.Item(19, 15).Value = 18000
$ws.Cells.Item(19, 16).Value = 18000
$ws.Cells.Item(19, 19).Value = 1000

# Row 20
$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, 3).Value = "Ñuble"
$ws.Cells.Item(20, 4).Value = 45001
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 5).Value = 16
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100107
$ws.Cells.Item(20, 8).Value = "Otros"
$ws.Cells.Item(20, 9).Value = 100107011
$ws.Cells.Item(20, 10).Value = "Tuna"
$ws.Cells.Item(20, 11).Value = "Sin especificar"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 60
$ws.Cells.Item(20, 14).Value = 17000
$ws.Cells.Item(20, 15).Value = 18000
$ws.Cells.Item(20, 16).Value = 17500
$ws.Cells.Item(20, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(20, 18).Value = "Región Metropolitana"
$ws.Cells.Item(20, 19).Value = 972
$ws.Cells.Item(20, 20).Value = 18
